$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.473.12"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.172.98"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "614.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.388"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.685"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.166.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.175"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -6.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.053.15"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.746.88"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.20"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.190.83"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.26"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  +36.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.361.44"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "74.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -8.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +36.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "529.93"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.90"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("E36").Value = "  -5.29%  "
$ws.Range("E37").Value = "  -10.43%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.27"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.82"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.126"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -9.90%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.374"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.28%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.44"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "172.30"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("E49").Value = "  -6.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.611"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.05"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.52%  "
